$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete row 11 ("Winkelwagen", 10, <blank>, "Nee") - entire row shift up
$ws.Rows.Item(11).Delete()

# 2. Fix values that changed in-place (after the shift):
#    row14 = MainComposer: B 0.2->0.5, C 0.2->0.5
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 0.5
#    row16 = Seeders and migrations: B 1.2->2, C 1.3->2
$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 2

# 3. Insert 11 new rows after row 17 ("Items"), before the filler/blank rows,
#    so the total row ends up at row 34 (10 new data rows + 1 extra blank filler row).
$ws.Range("A18:D28").EntireRow.Insert()

# Copy down the formatting (border style s="6") of row 17 into the freshly
# inserted rows so they match the rest of the table.
$ws.Range("A17:D17").Copy() | Out-Null
$ws.Range("A18:D28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4. Fill in the new rows (18-27) with the new work-breakdown items
$data = @(
  @("Product pagina", 2.5, 2.5, "Ja"),
  @("Toevoegen aan winkelwagen", 1.5, 1.5, "Ja"),
  @("Producten toevoegen aan winkelwagen", 1, 1, "Ja"),
  @("Product vooraad", 2, 2, "Ja"),
  @("Winkelwagen aantal en prijs weergave", 0.5, 0.5, "Ja"),
  @("Winkelwagen weergave", 2, 2, "Nee"),
  @("Winkelwagen Updaten", 1.5, $null, "Nee"),
  @("Winkelwagen verwijderen", 0.5, 0.5, "Ja"),
  @("Winkelwagen product verwijderen", 0.2, $null, "Nee"),
  @("Melding voor succesvol toevoegen van product", 0.2, 0.2, "Ja")
)

$r = 18
foreach ($row in $data) {
  $ws.Range("A$r").Value = $row[0]
  $ws.Range("B$r").Value = $row[1]
  if ($row[2] -ne $null) {
    $ws.Range("C$r").Value = $row[2]
  }
  $ws.Range("D$r").Value = $row[3]
  $r = $r + 1
}

# 5. Fix the Totaal row formulas (now at row 34): B34 sums the whole data block,
#    C34 only sums through row 22 (matches source workbook).
$ws.Range("B34").Formula = "=SUM(B3:B27)"
$ws.Range("C34").Formula = "=SUM(C3:C22)"

# 6. G1 formula already auto-adjusted to SUM(C34) by the row insert, but set
#    it explicitly to be safe.
$ws.Range("G1").Formula = "=SUM(C34)"

# 7. Column width / view adjustments
$ws.Columns.Item(1).ColumnWidth = 42.73697916666667
$ws.Columns.Item(13).ColumnWidth = 15.33984375

$ws.Range("C23").Select()
$excel.ActiveWindow.ScrollRow = 13
